$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.028.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.21%  '

$ws.Range("D3").Value = "'1.910.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.58%  '

$ws.Range("E4").Value = '  -0.70%  '

$ws.Range("D5").Value = "'316.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("E6").Value = '  -0.66%  '

$ws.Range("E8").Value = '  +0.35%  '

$ws.Range("D9").Value = "'0.07359"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.50%  '

$ws.Range("D10").Value = "'0.9347"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("D11").Value = "'20.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.38%  '

$ws.Range("D12").Value = "'0.07812"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").Value = "'1.912.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.29%  '

$ws.Range("D14").Value = "'5.507"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.40%  '

$ws.Range("D15").Value = "'6.641"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.46%  '

$ws.Range("D16").Value = "'92.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("D18").Value = "'0.000008879"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("D20").Value = "'28.056.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.98%  '

$ws.Range("D21").Value = "'14.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.73%  '

$ws.Range("D22").Value = "'5.171"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.20%  '

$ws.Range("D23").Value = "'2.146.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.43%  '

$ws.Range("E24").Value = '  +2.04%  '

$ws.Range("D25").Value = "'157.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").Value = "'1.911"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.67%  '

$ws.Range("D27").Value = "'18.51"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").Value = "'2.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.63%  '

$ws.Range("D29").Value = "'117.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.69%  '

$ws.Range("D30").Value = "'4.980"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.51%  '

$ws.Range("D31").Value = "'0.08956"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("E32").Value = '  -1.34%  '

$ws.Range("D33").Value = "'1.259"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.77%  '

$ws.Range("D34").Value = "'0.7750"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "'4.662"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.46%  '

$ws.Range("D36").Value = "'2.620"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.59%  '

$ws.Range("D37").Value = "'0.02050"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("E38").Value = '  -1.23%  '

$ws.Range("D39").Value = "'0.5539"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.26%  '

$ws.Range("D40").Value = "'0.05297"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.59%  '

$ws.Range("D41").Value = "'2.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.26%  '

$ws.Range("D42").Value = "'7.032"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("D43").Value = "'8.523"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.29%  '

$ws.Range("D45").Value = "'10.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("E46").Value = '  +5.43%  '

$ws.Range("D47").Value = "'0.4835"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.19%  '

$ws.Range("E48").Value = '  -0.69%  '

$ws.Range("D49").Value = "'1.652"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("D50").Value = "'68.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("D51").Value = "'0.06072"
$ws.Range("D51").Style = "Normal"
